# Apply the cibmtr-reporting-ig metadata update to the "Metadata" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value updates -------------------------------------------------
$ws.Range("B3").Value = "0.1.7"                                 # Version
$ws.Range("B6").Value = "draft"                                 # Status
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"             # Date
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"  # Contact
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"        # Contact (2nd row)

# --- insert a new "Jurisdiction" row after the Contact rows --------------
$ws.Rows.Item(12).Insert()

# copy the formatting of the row above so the new row matches the sheet's style
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
